{"js": "const replacements = [\n  [\"661\u00d78=5288\", \"364\u00d72=728\"],\n  [\"737\u00d74=2948\", \"117\u00d73=351\"],\n  [\"469\u00d79=4221\", \"505\u00d72=1010\"],\n  [\"790\u00d77=5530\", \"238\u00d77=1666\"],\n  [\"720\u00d78=5760\", \"958\u00d73=2874\"],\n  [\"763\u00d73=2289\", \"270\u00d72=540\"],\n  [\"964\u00d74=3856\", \"188\u00d78=1504\"],\n  [\"693\u00d78=5544\", \"156\u00d72=312\"],\n  [\"202\u00d72=404\", \"575\u00d74=2300\"],\n  [\"552\u00d76=3312\", \"385\u00d76=2310\"],\n  [\"590\u00d77=4130\", \"784\u00d73=2352\"],\n  [\"844\u00d73=2532\", \"206\u00d75=1030\"],\n  [\"551\u00d74=2204\", \"312\u00d74=1248\"],\n  [\"950\u00d77=6650\", \"756\u00d74=3024\"],\n  [\"145\u00d79=1305\", \"791\u00d77=5537\"],\n  [\"771\u00d73=2313\", \"239\u00d76=1434\"],\n  [\"217\u00d75=1085\", \"691\u00d74=2764\"],\n  [\"357\u00d75=1785\", \"501\u00d76=3006\"],\n  [\"345\u00d75=1725\", \"923\u00d74=3692\"],\n  [\"226\u00d73=678\", \"686\u00d73=2058\"],\n  [\"122\u00d75=610\", \"760\u00d76=4560\"],\n  [\"943\u00d73=2829\", \"447\u00d78=3576\"],\n  [\"139\u00d79=1251\", \"521\u00d72=1042\"],\n  [\"138\u00d73=414\", \"617\u00d77=4319\"],\n  [\"584\u00d76=3504\", \"976\u00d72=1952\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  ,@(\"661\u00d78=5288\", \"364\u00d72=728\")\n  ,@(\"737\u00d74=2948\", \"117\u00d73=351\")\n  ,@(\"469\u00d79=4221\", \"505\u00d72=1010\")\n  ,@(\"790\u00d77=5530\", \"238\u00d77=1666\")\n  ,@(\"720\u00d78=5760\", \"958\u00d73=2874\")\n  ,@(\"763\u00d73=2289\", \"270\u00d72=540\")\n  ,@(\"964\u00d74=3856\", \"188\u00d78=1504\")\n  ,@(\"693\u00d78=5544\", \"156\u00d72=312\")\n  ,@(\"202\u00d72=404\", \"575\u00d74=2300\")\n  ,@(\"552\u00d76=3312\", \"385\u00d76=2310\")\n  ,@(\"590\u00d77=4130\", \"784\u00d73=2352\")\n  ,@(\"844\u00d73=2532\", \"206\u00d75=1030\")\n  ,@(\"551\u00d74=2204\", \"312\u00d74=1248\")\n  ,@(\"950\u00d77=6650\", \"756\u00d74=3024\")\n  ,@(\"145\u00d79=1305\", \"791\u00d77=5537\")\n  ,@(\"771\u00d73=2313\", \"239\u00d76=1434\")\n  ,@(\"217\u00d75=1085\", \"691\u00d74=2764\")\n  ,@(\"357\u00d75=1785\", \"501\u00d76=3006\")\n  ,@(\"345\u00d75=1725\", \"923\u00d74=3692\")\n  ,@(\"226\u00d73=678\", \"686\u00d73=2058\")\n  ,@(\"122\u00d75=610\", \"760\u00d76=4560\")\n  ,@(\"943\u00d73=2829\", \"447\u00d78=3576\")\n  ,@(\"139\u00d79=1251\", \"521\u00d72=1042\")\n  ,@(\"138\u00d73=414\", \"617\u00d77=4319\")\n  ,@(\"584\u00d76=3504\", \"976\u00d72=1952\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}"}
